$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the note text shared across F20:F23 ("Channel recorded constant 63 for
# the entire day" -> "Channel recorded constant for the entire day but should
# not have"). All four cells shared the same string, so update all of them to
# keep them in sync.
$newText = "Channel recorded constant for the entire day but should not have"
$ws.Range("F20:F23").Value = $newText

# Reset the active sheet's selection back to the default (A1), clearing the
# previously-saved selection on F4.
$ws.Activate()
$ws.Range("A1").Select()
